{"js": "// Replace the division problems in the table with their new values.\n// Each original expression is unique within the document, so a simple\n// search-and-replace (matching the whole cell text) is safe here.\nconst replacements = [\n  [\"657\u00f76=\", \"236\u00f77=\"],\n  [\"798\u00f76=\", \"836\u00f75=\"],\n  [\"849\u00f72=\", \"193\u00f72=\"],\n  [\"397\u00f79=\", \"176\u00f78=\"],\n  [\"366\u00f76=\", \"395\u00f72=\"],\n  [\"866\u00f74=\", \"877\u00f74=\"],\n  [\"260\u00f72=\", \"543\u00f77=\"],\n  [\"748\u00f74=\", \"728\u00f75=\"],\n  [\"379\u00f72=\", \"671\u00f74=\"],\n  [\"325\u00f74=\", \"113\u00f76=\"],\n  [\"671\u00f76=\", \"199\u00f73=\"],\n  [\"171\u00f79=\", \"529\u00f78=\"],\n  [\"879\u00f72=\", \"477\u00f77=\"],\n  [\"506\u00f77=\", \"350\u00f75=\"],\n  [\"173\u00f77=\", \"420\u00f79=\"],\n  [\"278\u00f75=\", \"192\u00f79=\"],\n  [\"873\u00f75=\", \"386\u00f76=\"],\n  [\"836\u00f72=\", \"917\u00f79=\"],\n  [\"930\u00f73=\", \"281\u00f75=\"],\n  [\"725\u00f78=\", \"446\u00f79=\"],\n  [\"347\u00f75=\", \"234\u00f73=\"],\n  [\"168\u00f73=\", \"531\u00f73=\"],\n  [\"790\u00f79=\", \"153\u00f79=\"],\n  [\"241\u00f73=\", \"171\u00f78=\"],\n  [\"731\u00f77=\", \"361\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with their new values.\n# Each original expression is unique within the document, so a simple\n# Find/Replace (ReplaceAll) on the whole document content is safe here.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"657\u00f76=\"; New = \"236\u00f77=\" },\n    @{ Old = \"798\u00f76=\"; New = \"836\u00f75=\" },\n    @{ Old = \"849\u00f72=\"; New = \"193\u00f72=\" },\n    @{ Old = \"397\u00f79=\"; New = \"176\u00f78=\" },\n    @{ Old = \"366\u00f76=\"; New = \"395\u00f72=\" },\n    @{ Old = \"866\u00f74=\"; New = \"877\u00f74=\" },\n    @{ Old = \"260\u00f72=\"; New = \"543\u00f77=\" },\n    @{ Old = \"748\u00f74=\"; New = \"728\u00f75=\" },\n    @{ Old = \"379\u00f72=\"; New = \"671\u00f74=\" },\n    @{ Old = \"325\u00f74=\"; New = \"113\u00f76=\" },\n    @{ Old = \"671\u00f76=\"; New = \"199\u00f73=\" },\n    @{ Old = \"171\u00f79=\"; New = \"529\u00f78=\" },\n    @{ Old = \"879\u00f72=\"; New = \"477\u00f77=\" },\n    @{ Old = \"506\u00f77=\"; New = \"350\u00f75=\" },\n    @{ Old = \"173\u00f77=\"; New = \"420\u00f79=\" },\n    @{ Old = \"278\u00f75=\"; New = \"192\u00f79=\" },\n    @{ Old = \"873\u00f75=\"; New = \"386\u00f76=\" },\n    @{ Old = \"836\u00f72=\"; New = \"917\u00f79=\" },\n    @{ Old = \"930\u00f73=\"; New = \"281\u00f75=\" },\n    @{ Old = \"725\u00f78=\"; New = \"446\u00f79=\" },\n    @{ Old = \"347\u00f75=\"; New = \"234\u00f73=\" },\n    @{ Old = \"168\u00f73=\"; New = \"531\u00f73=\" },\n    @{ Old = \"790\u00f79=\"; New = \"153\u00f79=\" },\n    @{ Old = \"241\u00f73=\"; New = \"171\u00f78=\" },\n    @{ Old = \"731\u00f77=\"; New = \"361\u00f75=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
